# "Se actualiza puntos trabajadores"
#
# Updates the points totals on the two "current month" sheets and leaves the
# active-cell selection where the author last clicked before saving.

$wb = $excel.ActiveWorkbook

# --- Sheet "Mes actual cel": points column (C) for the first three rows ---
$wsCel = $wb.Worksheets.Item("Mes actual cel")
$wsCel.Range("C2").Value = 63.865
$wsCel.Range("C3").Value = 63.865
$wsCel.Range("C4").Value = 63.865

# --- Sheet "Mes actual emp": points column (C) for the single data row ---
$wsEmp = $wb.Worksheets.Item("Mes actual emp")
$wsEmp.Range("C2").Value = 141.89

# Leave the cursor on C2 of "Mes actual emp" (matches the saved selection),
# then return to "Mes actual cel" and park the cursor on C4 so that sheet
# stays the active/selected tab, matching the workbook's final state.
$wsEmp.Range("C2").Select()

$wsCel.Activate()
$wsCel.Range("C4").Select()
